
# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# Column BF ("Date") held the literal text "6-26-2011-12" on every data
# row; correct it to the proper ISO date string "2012-06-26" while
# keeping the cells as plain text (not an Excel date serial) and without
# disturbing their existing (unstyled) formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BF2:BF31 holds the 30 data rows with the bad date text.
$dateRange = $ws.Range("BF2:BF31")

# Pre-format as Text so Excel doesn't reinterpret the ISO-looking string
# "2012-06-26" as a date serial when we assign it.
$dateRange.NumberFormat = "@"
$dateRange.Value = "2012-06-26"

# The original cells carried no explicit style (default "Normal"); put
# that back now that the text-number-format detour is done, so the only
# observable change is the corrected cell text.
$dateRange.Style = "Normal"
